$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Grab the percentage number format already used by column E so the new
# column F cells reuse the existing style (xf index 2) instead of Excel
# creating a brand-new cell style entry.
$pctFormat = $ws.Range("E2").NumberFormat

# New "day" column F of percentages, mirroring column E's layout.
$ws.Range("F2").Value = 0.69399999999999995
$ws.Range("F2").NumberFormat = $pctFormat

$ws.Range("F3").Value = 0.112
$ws.Range("F3").NumberFormat = $pctFormat

$ws.Range("F4").Value = 0.46600000000000003
$ws.Range("F4").NumberFormat = $pctFormat

$ws.Range("F5").Value = 0.17199999999999999
$ws.Range("F5").NumberFormat = $pctFormat

$ws.Range("F6").Value = 0.39100000000000001
$ws.Range("F6").NumberFormat = $pctFormat

$ws.Range("F8").Value = 0.44700000000000001
$ws.Range("F8").NumberFormat = $pctFormat

# A little note from LilyAnne, added as a new row underneath the table.
$ws.Range("F9").Value = "^I feel like this is wrong, but it's what my computer says. -LilyAnne"

# Move the active selection to where it ended up after the edits.
$ws.Range("F9").Select()
